$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 23
$ws.Range("B5").Value = 36
$ws.Range("B6").Select()
